# Commit: "added back response curve functionality"
#
# 1) Column B (YEAR_WEEK) labels were using placeholder year tokens
#    Y1 / Y2 / Y3 — restore the real calendar years: 2020 / 2021 / 2022.
# 2) Column H (ROS) had literal "inf" text left over in some rows from a
#    previous (broken) run of the response-curve calc; now that the
#    response curve logic is back, those rows resolve to a real number (0).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$yearMap = @{ "Y1" = "2020"; "Y2" = "2021"; "Y3" = "2022" }

$dims = $ws.UsedRange
$lastRow = $dims.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {

    # --- Column B: Y{n}_{week} -> {year}_{week} ---
    $bCell = $ws.Cells.Item($r, 2)
    $bVal = $bCell.Value2
    if ($null -ne $bVal) {
        $bText = [string]$bVal
        if ($bText -match "^(Y[123])_(\d+)$") {
            $prefix = $matches[1]
            $week = $matches[2]
            $bCell.Value = $yearMap[$prefix] + "_" + $week
        }
    }

    # --- Column H: literal "inf" text -> numeric 0 ---
    $hCell = $ws.Cells.Item($r, 8)
    $hVal = $hCell.Value2
    if ($null -ne $hVal) {
        if ($hVal -is [string]) {
            if ($hVal -eq "inf") {
                $hCell.Value = 0
            }
        }
    }
}
